$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '25.506.55'
Set-TextValue 'E2' '  +1.65%  '
Set-TextValue 'D3' '1.663.32'
Set-TextValue 'E3' '  +0.76%  '
Set-TextValue 'D4' '0.9991'
Set-TextValue 'E4' '  -0.11%  '
Set-TextValue 'D5' '236.92'
Set-TextValue 'E5' '  -0.24%  '
Set-TextValue 'E6' '  -0.03%  '
Set-TextValue 'E7' '  +0.24%  '
Set-TextValue 'D8' '0.2623'
Set-TextValue 'E8' '  +0.34%  '
Set-TextValue 'D9' '0.06156'
Set-TextValue 'E9' '  +2.50%  '
Set-TextValue 'D10' '0.07103'
Set-TextValue 'E10' '  -1.19%  '
Set-TextValue 'D11' '1.662.89'
Set-TextValue 'E11' '  +0.65%  '
Set-TextValue 'E12' '  -0.91%  '
Set-TextValue 'D13' '0.5867'
Set-TextValue 'E13' '  -5.97%  '
Set-TextValue 'D14' '4.368'
Set-TextValue 'E14' '  -5.17%  '
Set-TextValue 'D15' '74.72'
Set-TextValue 'E15' '  +1.79%  '
Set-TextValue 'E16' '  +0.06%  '
Set-TextValue 'D17' '0.9999'
Set-TextValue 'E17' '  -0.01%  '
Set-TextValue 'D18' '25.492.82'
Set-TextValue 'E18' '  +1.63%  '
Set-TextValue 'D19' '0.000006738'
Set-TextValue 'E19' '  +1.96%  '
Set-TextValue 'D20' '11.46'
Set-TextValue 'D21' '1.873.34'
Set-TextValue 'E21' '  +0.66%  '
Set-TextValue 'D22' '4.426'
Set-TextValue 'E22' '  -1.08%  '
Set-TextValue 'D23' '8.668'
Set-TextValue 'E23' '  +0.67%  '
Set-TextValue 'D24' '5.269'
Set-TextValue 'E24' '  -0.41%  '
Set-TextValue 'D25' '133.59'
Set-TextValue 'E25' '  +0.34%  '
Set-TextValue 'D26' '15.05'
Set-TextValue 'E26' '  +0.62%  '
Set-TextValue 'D27' '1.387'
Set-TextValue 'E27' '  -0.56%  '
Set-TextValue 'D28' '105.31'
Set-TextValue 'E28' '  +1.76%  '
Set-TextValue 'D29' '1.712'
Set-TextValue 'E29' '  +1.85%  '
Set-TextValue 'D30' '3.945'
Set-TextValue 'E30' '  +4.55%  '
Set-TextValue 'D31' '3.665'
Set-TextValue 'E31' '  +2.66%  '
Set-TextValue 'D32' '0.07661'
Set-TextValue 'E32' '  -3.28%  '
Set-TextValue 'D33' '0.9993'
Set-TextValue 'E33' '  -0.03%  '
Set-TextValue 'D34' '0.04211'
Set-TextValue 'E34' '  -8.37%  '
Set-TextValue 'E35' '  +0.57%  '
Set-TextValue 'D36' '0.6098'
Set-TextValue 'E36' '  +5.67%  '
Set-TextValue 'D37' '0.9500'
Set-TextValue 'E37' '  +0.52%  '
Set-TextValue 'E38' '  -0.57%  '
Set-TextValue 'D39' '0.8665'
Set-TextValue 'E39' '  +3.68%  '
Set-TextValue 'D40' '0.9996'
Set-TextValue 'E40' '  -0.09%  '
Set-TextValue 'D41' '1.853'
Set-TextValue 'E41' '  +1.27%  '
Set-TextValue 'E42' '  -5.54%  '
Set-TextValue 'D43' '96.84'
Set-TextValue 'E43' '  -2.29%  '
Set-TextValue 'D44' '0.3758'
Set-TextValue 'E44' '  +1.02%  '
Set-TextValue 'D45' '4.760'
Set-TextValue 'E45' '  -1.15%  '
Set-TextValue 'D46' '0.1126'
Set-TextValue 'E46' '  -1.14%  '
Set-TextValue 'D47' '6.205'
Set-TextValue 'E47' '  +1.54%  '
Set-TextValue 'E48' '  +1.33%  '
Set-TextValue 'D49' '29.66'
Set-TextValue 'E49' '  -0.53%  '
Set-TextValue 'E50' '  -0.05%  '
Set-TextValue 'D51' '0.9992'
Set-TextValue 'E51' '  -0.08%  '
